# Season record columns: Wins / Losses / Ties
# Adds AD/AE/AF headers (copying the existing header formatting from A1)
# and fills the season record (77 wins, 85 losses, 0 ties) for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 44

# --- Headers (row 1) ---------------------------------------------------
# Copy the header cell's formatting (bold, centered, thin border) onto the
# new header cells, then set their text.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2..44) ---------------------------------------------------
for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 30).Value = 77   # AD -> Wins
    $ws.Cells.Item($row, 31).Value = 85   # AE -> Losses
    $ws.Cells.Item($row, 32).Value = 0    # AF -> Ties
}
